$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column B for rows 1-9 with "D", except row 8 which gets "Logic done"
# (this mirrors marking each algorithm's status, with row 8 -
#  "Check weather a given number is positive or negative" - marked specially)
for ($r = 1; $r -le 9; $r++) {
    if ($r -eq 8) {
        $ws.Cells.Item($r, 2).Value = "Logic done"
    } else {
        $ws.Cells.Item($r, 2).Value = "D"
    }
}

# Move the active selection to F9, as left after the edits
$ws.Range("F9").Select()

# Set the page to portrait orientation (page setup touched)
$ws.PageSetup.Orientation = 1
